# Attendance Marking Feature
# Populates the attendance log with check-in/check-out entries and
# formats the header row + timestamp columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row formatting: bold, boxed border, centered/top-aligned ---
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop

# --- Attendance data rows ---
$names = @(
    "Harshit Saraswat",
    "Harshit Saraswat",
    "Harshit Saraswat",
    "Harshit Saraswat",
    "Harshit Saraswat"
)
$checkIn = @(
    44513.82614364583,
    44513.83076476852,
    44513.83081765047,
    44513.83127752315,
    44513.83187018075
)
$checkOut = @(
    44513.83074774306,
    44513.83080909722,
    44513.83084804398,
    44513.8313628125,
    $null
)

for ($i = 0; $i -lt 5; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $names[$i]

    $ws.Cells.Item($r, 2).Value = $checkIn[$i]
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    if ($checkOut[$i] -ne $null) {
        $ws.Cells.Item($r, 3).Value = $checkOut[$i]
        $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
}

# --- Page margins (inches -> points: 1in = 72pt) ---
$ws.PageSetup.LeftMargin = 54      # 0.75"
$ws.PageSetup.RightMargin = 54     # 0.75"
$ws.PageSetup.TopMargin = 72       # 1"
$ws.PageSetup.BottomMargin = 72    # 1"
$ws.PageSetup.HeaderMargin = 36    # 0.5"
$ws.PageSetup.FooterMargin = 36    # 0.5"

# --- Reset selection to A1 (matches the freshly generated sheet state) ---
[void]$ws.Range("A1").Select()
